# Commit: "Using internal column widths in pptx writer tables (#9392)"
#
# The pptx table writer used to always divide the graphic-frame's total
# width evenly across every column. Now it uses the individual column
# widths supplied when the table was authored (falling back to the even
# split only when no widths are available). For the table on this slide
# that means both of its columns grow from 197pt (2,501,900 EMU) to
# 198pt (2,514,600 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the (first) shape on the slide that holds a table.
$tblShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tblShape = $candidate
        break
    }
}

$tbl = $tblShape.Table

for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
    $tbl.Columns.Item($c).Width = 198
}
